$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists accounts ordered by Saldo descending. Two rows (AHMAD and
# GUILHERME) need to move up, ahead of the LEVI row, and LEVI's balance
# value changes from 100064.95 to 64000. Also the MIRELLA / 001651617 row
# is removed entirely.

# Step 1: Insert two blank rows above row 3 (currently holding LEVI), to
# make room for the AHMAD and GUILHERME rows that move above LEVI.
$ws.Rows("3:4").Insert()

# Step 2: Populate the new rows 3 and 4 with the AHMAD / GUILHERME data.
# Force text format on the account-number column so the leading zeros in
# the account numbers are preserved (matching the original inline-string
# cells).
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "004368468"
$ws.Cells.Item(3, 2).Value = "AHMAD"
$ws.Cells.Item(3, 3).Value = 88330.09

$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "004948033"
$ws.Cells.Item(4, 2).Value = "GUILHERME"
$ws.Cells.Item(4, 3).Value = 84855.32

# Step 3: The original LEVI row has shifted down to row 5; update its
# Saldo value to the new amount.
$ws.Cells.Item(5, 3).Value = 64000

# Step 4: The original AHMAD row (now row 6) and GUILHERME row (now row 7)
# are duplicates of the ones we just inserted above LEVI; remove them.
$ws.Rows("6:7").Delete()

# Step 5: Remove the MIRELLA / 001651617 row entirely (still at row 16,
# since steps 1-4 net to zero row-count change above it).
$ws.Rows(16).Delete()
